$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.2890784740448
$ws.Range("B1").Value = 2.68403434753418
$ws.Range("C1").Value = 2.384029865264893
$ws.Range("D1").Value = 2.616694450378418
$ws.Range("E1").Value = 3.253950595855713
